$wb = $excel.ActiveWorkbook

# Delete the empty "Sheet1" worksheet (keeps the data worksheet "Sheet2")
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Sheet1").Delete()

# Rename the remaining data worksheet to "Sheet1"
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "Sheet1"

# Update the selection on the remaining sheet
$ws.Activate()
$ws.Range("E25").Select()
